# Add Presentation Chunk to awesomeCV_amb
# - adds a new "Education" worksheet (first tab) with degree history
# - populates the pre-existing (empty) "Presentations" worksheet and
#   relocates it to sit right after "Teaching"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "Education" sheet, inserted before "Experience" (-> first tab)
# ---------------------------------------------------------------------
$experience = $wb.Worksheets.Item("Experience")
$edu = $wb.Worksheets.Add($experience)
$edu.Name = "Education"

$edu.Range("A1").Value = "role-title"
$edu.Range("B1").Value = "department-name"
$edu.Range("C1").Value = "end-date.year.value"
$edu.Range("D1").Value = "organization.name"
$edu.Range("E1").Value = "organization.address.city"
$edu.Range("F1").Value = "organization.address.region"
$edu.Range("G1").Value = "organization.address.country"
$edu.Range("H1").Value = "thesis"

$edu.Range("A2").Value = "Ph.D"
$edu.Range("B2").Value = "Neuroscience"
$edu.Range("C2").Value = "In Progress"
$edu.Range("D2").Value = "École Polytechnique Fédérale de Lausanne"
$edu.Range("E2").Value = "Lausanne"
$edu.Range("G2").Value = "Switerland"
$edu.Range("H2").Value = "Testing the theory of Epigenetic Priming in Fear Memory Conditioning"

$edu.Range("A3").Value = "M.Sc"
$edu.Range("B3").Value = "Bioinformatics"
$edu.Range("C3").Value = 2013
$edu.Range("D3").Value = "University of Oregon"
$edu.Range("E3").Value = "Eugene"
$edu.Range("F3").Value = "Oregon"
$edu.Range("G3").Value = "USA"
$edu.Range("H3").Value = "Effects of the splicing inhibitor, Isoginkgetin, on human Telomerase RNA"

$edu.Range("A4").Value = "B.Sc"
$edu.Range("B4").Value = "Biology, Minor in Chemistry"
$edu.Range("C4").Value = 2012
$edu.Range("D4").Value = "University of Oregon"
$edu.Range("E4").Value = "Eugene"
$edu.Range("F4").Value = "Oregon"
$edu.Range("G4").Value = "USA"

$edu.Columns.Item(1).ColumnWidth = 10.1640625
$edu.Range("B1:D1").ColumnWidth = 20.83203125
$edu.Columns.Item(5).ColumnWidth = 21.5
$edu.Columns.Item(6).ColumnWidth = 24.33203125
$edu.Columns.Item(7).ColumnWidth = 23.6640625

$edu.Range("F2").Select()

# ---------------------------------------------------------------------
# 2. Fill in the (already existing, but empty) "Presentations" sheet
# ---------------------------------------------------------------------
$pres = $wb.Worksheets.Item("Presentations")

$pres.Range("A1").Value = "Symposium"
$pres.Range("B1").Value = "Date"
$pres.Range("C1").Value = "Location"

$pres.Range("A2").Value = "MCCS 2019"
$pres.Range("B2").Value = "Oct. 18, 2019"
$pres.Range("C2").Value = "Chicago, Illinois"

$pres.Columns.Item(2).ColumnWidth = 15.5
$pres.Columns.Item(3).ColumnWidth = 14

# Move "Presentations" so it sits right after "Teaching"
$teaching = $wb.Worksheets.Item("Teaching")
$pres.Move([System.Reflection.Missing]::Value, $teaching)

# Re-resolve by name: after Move(), old object handles track the
# *position* they used to occupy rather than the sheet identity.
$pres = $wb.Worksheets.Item("Presentations")
$pres.Activate()
$pres.Range("B2").Select()
